$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the values in row 5 (B5:AH5) to 2 decimal places ("custom accuracy").
# Columns M5, AB5, AC5 already have <=2 decimals so they remain unchanged.
$ws.Range("B5").Value = 18.74
$ws.Range("C5").Value = 13.73
$ws.Range("D5").Value = 1.18
$ws.Range("E5").Value = 40.71
$ws.Range("F5").Value = 33.26
$ws.Range("G5").Value = 14.75
$ws.Range("H5").Value = 58.36
$ws.Range("I5").Value = 22.69
$ws.Range("J5").Value = 10.04
$ws.Range("K5").Value = 14.82
$ws.Range("L5").Value = 16.34
$ws.Range("M5").Value = 17.2
$ws.Range("N5").Value = 4.71
$ws.Range("O5").Value = 14.66
$ws.Range("P5").Value = 20.84
$ws.Range("Q5").Value = 12.42
$ws.Range("R5").Value = 0.83
$ws.Range("S5").Value = 0.79
$ws.Range("T5").Value = 215.87
$ws.Range("U5").Value = 41.04
$ws.Range("V5").Value = 13.53
$ws.Range("W5").Value = 27.51
$ws.Range("X5").Value = 14.4
$ws.Range("Y5").Value = 2.24
$ws.Range("Z5").Value = 28.15
$ws.Range("AA5").Value = 11.95
$ws.Range("AB5").Value = 10.64
$ws.Range("AC5").Value = 12.5
$ws.Range("AD5").Value = 17.11
$ws.Range("AE5").Value = 0.56
$ws.Range("AF5").Value = 52.98
$ws.Range("AG5").Value = 7.59
$ws.Range("AH5").Value = 16.92

# Remove the last data row (row 6) entirely; this also shrinks the used
# range / dimension from A1:AH6 down to A1:AH5.
$ws.Rows.Item(6).Delete()
